# Apply the artfynd.xlsx corrections: reorder a handful of observation
# records (rows 8-39 area) back to their correct row positions, and bump
# the taxon sort-order code (column B) for Garnlav / Rosenticka / Ullticka
# records that were using a stale "Taxonsorteringsordning" value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 131244272
$ws.Range("B8").Value = 79245
$ws.Range("Q8").Value = 613374
$ws.Range("R8").Value = 6998037
$ws.Range("S8").Value = 10
$ws.Range("AC8").Value = $null

# Row 9
$ws.Range("A9").Value = 131244269
$ws.Range("B9").Value = 79245
$ws.Range("Q9").Value = 613333
$ws.Range("R9").Value = 6998073
$ws.Range("S9").Value = 25
$ws.Range("AC9").Value = 'Rikligt med garnlav inom ett område på 25m'

# Row 10
$ws.Range("A10").Value = 131244300
$ws.Range("B10").Value = 79245
$ws.Range("Q10").Value = 613444
$ws.Range("R10").Value = 6998046
$ws.Range("AJ10").Value = 'tall'
$ws.Range("AK10").Value = 'Pinus sylvestris'
$ws.Range("AO10").Value = 'Pinus sylvestris'

# Row 11
$ws.Range("A11").Value = 131244279
$ws.Range("B11").Value = 79245
$ws.Range("Q11").Value = 613427
$ws.Range("R11").Value = 6998062
$ws.Range("AJ11").Value = 'gran'
$ws.Range("AK11").Value = 'Picea abies'
$ws.Range("AO11").Value = 'Picea abies'

# Row 13
$ws.Range("B13").Value = 79245

# Row 14
$ws.Range("B14").Value = 79245

# Row 15
$ws.Range("A15").Value = 131244281
$ws.Range("B15").Value = 79245
$ws.Range("Q15").Value = 613346
$ws.Range("R15").Value = 6998128
$ws.Range("S15").Value = 50
$ws.Range("AC15").Value = 'Rikliga mängder garnlav på främst tall men även gran inom ett område på ca 50 m.'
$ws.Range("AJ15").Value = $null
$ws.Range("AK15").Value = $null
$ws.Range("AO15").Value = $null

# Row 16
$ws.Range("A16").Value = 131244280
$ws.Range("B16").Value = 79245
$ws.Range("Q16").Value = 613405
$ws.Range("R16").Value = 6998110
$ws.Range("S16").Value = 10
$ws.Range("AC16").Value = $null
$ws.Range("AJ16").Value = 'gran'
$ws.Range("AK16").Value = 'Picea abies'
$ws.Range("AO16").Value = 'Picea abies'

# Row 17
$ws.Range("A17").Value = 131244253
$ws.Range("B17").Value = 57884
$ws.Range("E17").Value = 100109
$ws.Range("F17").Value = 'Tretåig hackspett'
$ws.Range("G17").Value = 'Picoides tridactylus'
$ws.Range("H17").Value = '(Linnaeus, 1758)'
$ws.Range("M17").Value = 'färska spår'
$ws.Range("Q17").Value = 613354
$ws.Range("R17").Value = 6998128
$ws.Range("AC17").Value = 'Färska ringhack på tall'
$ws.Range("AJ17").Value = $null
$ws.Range("AK17").Value = $null
$ws.Range("AO17").Value = $null

# Row 18
$ws.Range("B18").Value = 79245

# Row 19
$ws.Range("A19").Value = 131244274
$ws.Range("B19").Value = 79245
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = 'Garnlav'
$ws.Range("G19").Value = 'Alectoria sarmentosa'
$ws.Range("H19").Value = '(Ach.) Ach.'
$ws.Range("M19").Value = $null
$ws.Range("Q19").Value = 613467
$ws.Range("R19").Value = 6997940
$ws.Range("AC19").Value = $null
$ws.Range("AJ19").Value = 'tall'
$ws.Range("AK19").Value = 'Pinus sylvestris'
$ws.Range("AO19").Value = 'Pinus sylvestris'

# Row 20
$ws.Range("A20").Value = 131244257
$ws.Range("B20").Value = 57881
$ws.Range("E20").Value = 100049
$ws.Range("F20").Value = 'Spillkråka'
$ws.Range("G20").Value = 'Dryocopus martius'
$ws.Range("H20").Value = '(Linnaeus, 1758)'
$ws.Range("M20").Value = 'färska spår'
$ws.Range("Q20").Value = 613343
$ws.Range("R20").Value = 6998121
$ws.Range("AJ20").Value = $null
$ws.Range("AK20").Value = $null
$ws.Range("AO20").Value = $null

# Row 21
$ws.Range("A21").Value = 131244289
$ws.Range("B21").Value = 79245
$ws.Range("E21").Value = 6425
$ws.Range("F21").Value = 'Garnlav'
$ws.Range("G21").Value = 'Alectoria sarmentosa'
$ws.Range("H21").Value = '(Ach.) Ach.'
$ws.Range("M21").Value = $null
$ws.Range("Q21").Value = 613338
$ws.Range("R21").Value = 6998224
$ws.Range("AJ21").Value = 'tall'
$ws.Range("AK21").Value = 'Pinus sylvestris'
$ws.Range("AO21").Value = 'Pinus sylvestris'

# Row 22
$ws.Range("A22").Value = 131244278
$ws.Range("B22").Value = 79245
$ws.Range("Q22").Value = 613406
$ws.Range("R22").Value = 6998050
$ws.Range("S22").Value = 10
$ws.Range("AC22").Value = $null
$ws.Range("AJ22").Value = 'gran'
$ws.Range("AK22").Value = 'Picea abies'
$ws.Range("AO22").Value = 'Picea abies'

# Row 23
$ws.Range("A23").Value = 131244276
$ws.Range("B23").Value = 79245
$ws.Range("Q23").Value = 613444
$ws.Range("R23").Value = 6998014
$ws.Range("S23").Value = 25
$ws.Range("AC23").Value = 'Rilkigt med garnlav inom ett ca 25m område'
$ws.Range("AJ23").Value = 'tall'
$ws.Range("AK23").Value = 'Pinus sylvestris'
$ws.Range("AO23").Value = 'Pinus sylvestris'

# Row 24
$ws.Range("A24").Value = 131244288
$ws.Range("B24").Value = 79245
$ws.Range("Q24").Value = 613331
$ws.Range("R24").Value = 6998221

# Row 25
$ws.Range("A25").Value = 131244290
$ws.Range("B25").Value = 79245
$ws.Range("Q25").Value = 613327
$ws.Range("R25").Value = 6998224

# Row 26
$ws.Range("B26").Value = 79245

# Row 28
$ws.Range("B28").Value = 92108

# Row 31
$ws.Range("B31").Value = 79245

# Row 32
$ws.Range("A32").Value = 131244263
$ws.Range("B32").Value = 79245
$ws.Range("Q32").Value = 613379
$ws.Range("R32").Value = 6998218
$ws.Range("S32").Value = 25
$ws.Range("AC32").Value = 'Rikliga mängder garnlav på gran inom ett område på 25m.'
$ws.Range("AJ32").Value = 'gran'
$ws.Range("AK32").Value = 'Picea abies'
$ws.Range("AO32").Value = 'Picea abies'

# Row 33
$ws.Range("A33").Value = 131244266
$ws.Range("B33").Value = 79245
$ws.Range("Q33").Value = 613400
$ws.Range("R33").Value = 6997964
$ws.Range("S33").Value = 10
$ws.Range("AC33").Value = $null
$ws.Range("AJ33").Value = 'tall'
$ws.Range("AK33").Value = 'Pinus sylvestris'
$ws.Range("AO33").Value = 'Pinus sylvestris'

# Row 34
$ws.Range("B34").Value = 79245

# Row 35
$ws.Range("B35").Value = 91810

# Row 36
$ws.Range("B36").Value = 79245

# Row 37
$ws.Range("A37").Value = 131244250
$ws.Range("B37").Value = 57884
$ws.Range("E37").Value = 100109
$ws.Range("F37").Value = 'Tretåig hackspett'
$ws.Range("G37").Value = 'Picoides tridactylus'
$ws.Range("H37").Value = '(Linnaeus, 1758)'
$ws.Range("M37").Value = 'färska spår'
$ws.Range("Q37").Value = 613387
$ws.Range("R37").Value = 6998025
$ws.Range("AC37").Value = 'Färska ringhack på tall'
$ws.Range("AJ37").Value = $null
$ws.Range("AK37").Value = $null
$ws.Range("AO37").Value = $null

# Row 38
$ws.Range("A38").Value = 131244297
$ws.Range("B38").Value = 79245
$ws.Range("E38").Value = 6425
$ws.Range("F38").Value = 'Garnlav'
$ws.Range("G38").Value = 'Alectoria sarmentosa'
$ws.Range("H38").Value = '(Ach.) Ach.'
$ws.Range("M38").Value = $null
$ws.Range("Q38").Value = 613464
$ws.Range("R38").Value = 6998117
$ws.Range("AC38").Value = $null
$ws.Range("AJ38").Value = 'tall'
$ws.Range("AK38").Value = 'Pinus sylvestris'
$ws.Range("AO38").Value = 'Pinus sylvestris'

# Row 39
$ws.Range("B39").Value = 79245
